$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$s.Shapes.Item(3).Table.ApplyStyle("{CE348E7A-7A39-420F-B0B8-1D0B22042551}")
